# Auto-applied numeric updates reproducing the upstream commit
# "chore: update Sheets via scheduled runner" (Zodiark_Profits workbook).
$wb = $excel.ActiveWorkbook

$changes = @(
    @{Sheet="ALC"; Cell="H19"; Value=2173.1667}
    @{Sheet="ALC"; Cell="J19"; Value=2598.625}
    @{Sheet="ALC"; Cell="L19"; Value=2598.625}
    @{Sheet="ALC"; Cell="N19"; Value=-2948.625}
    @{Sheet="ALC"; Cell="H40"; Value=7367.2856}
    @{Sheet="ALC"; Cell="J40"; Value=10830.25}
    @{Sheet="ALC"; Cell="L40"; Value=10830.25}
    @{Sheet="ALC"; Cell="N40"; Value=-11180.25}
    @{Sheet="ALC"; Cell="H116"; Value=3891.7}
    @{Sheet="ALC"; Cell="I116"; Value=3792.8}
    @{Sheet="ALC"; Cell="J116"; Value=3990.6}
    @{Sheet="ALC"; Cell="K116"; Value=3792.8}
    @{Sheet="ALC"; Cell="L116"; Value=3990.6}
    @{Sheet="ALC"; Cell="M116"; Value=-350.8000000000002}
    @{Sheet="ALC"; Cell="N116"; Value=-10874.6}
    @{Sheet="ALC"; Cell="H124"; Value=89999.5}
    @{Sheet="ALC"; Cell="J124"; Value=89999.5}
    @{Sheet="ALC"; Cell="L124"; Value=89999.5}
    @{Sheet="ALC"; Cell="N124"; Value=-99819.5}
    @{Sheet="ALC"; Cell="H126"; Value=99956.836}
    @{Sheet="ALC"; Cell="J126"; Value=99956.836}
    @{Sheet="ALC"; Cell="L126"; Value=99956.836}
    @{Sheet="ALC"; Cell="N126"; Value=-109836.836}
    @{Sheet="ALC"; Cell="H132"; Value=5121.1177}
    @{Sheet="ALC"; Cell="I132"; Value=3234.7144}
    @{Sheet="ALC"; Cell="J132"; Value=13924.333}
    @{Sheet="ALC"; Cell="K132"; Value=9704.143199999999}
    @{Sheet="ALC"; Cell="L132"; Value=41772.999}
    @{Sheet="ALC"; Cell="M132"; Value=-7174.143199999999}
    @{Sheet="ALC"; Cell="N132"; Value=-46832.999}
    @{Sheet="ALC"; Cell="H137"; Value=1680.5294}
    @{Sheet="ALC"; Cell="I137"; Value=1412.4}
    @{Sheet="ALC"; Cell="J137"; Value=1892.2106}
    @{Sheet="ALC"; Cell="K137"; Value=4237.200000000001}
    @{Sheet="ALC"; Cell="L137"; Value=5676.6318}
    @{Sheet="ALC"; Cell="M137"; Value=-1687.200000000001}
    @{Sheet="ALC"; Cell="N137"; Value=-10776.6318}
    @{Sheet="ARM"; Cell="H2"; Value=2113.9}
    @{Sheet="ARM"; Cell="I2"; Value=1172.5264}
    @{Sheet="ARM"; Cell="K2"; Value=1172.5264}
    @{Sheet="ARM"; Cell="M2"; Value=-1059.5264}
    @{Sheet="ARM"; Cell="H32"; Value=1941}
    @{Sheet="ARM"; Cell="I32"; Value=1578.6136}
    @{Sheet="ARM"; Cell="K32"; Value=1578.6136}
    @{Sheet="ARM"; Cell="M32"; Value=-1291.6136}
    @{Sheet="ARM"; Cell="H45"; Value=2499.5}
    @{Sheet="ARM"; Cell="I45"; Value=2666}
    @{Sheet="ARM"; Cell="K45"; Value=2666}
    @{Sheet="ARM"; Cell="M45"; Value=-2289}
    @{Sheet="ARM"; Cell="H63"; Value=1985.25}
    @{Sheet="ARM"; Cell="I63"; Value=1985.25}
    @{Sheet="ARM"; Cell="K63"; Value=1985.25}
    @{Sheet="ARM"; Cell="M63"; Value=-1299.25}
    @{Sheet="ARM"; Cell="H66"; Value=1985.25}
    @{Sheet="ARM"; Cell="I66"; Value=1985.25}
    @{Sheet="ARM"; Cell="K66"; Value=9926.25}
    @{Sheet="ARM"; Cell="M66"; Value=-6494.25}
    @{Sheet="ARM"; Cell="H116"; Value=2113.9}
    @{Sheet="ARM"; Cell="I116"; Value=1172.5264}
    @{Sheet="ARM"; Cell="K116"; Value=1172.5264}
    @{Sheet="ARM"; Cell="M116"; Value=1121.4736}
    @{Sheet="ARM"; Cell="L118"; Value=0}
    @{Sheet="ARM"; Cell="H118"; Value=0}
    @{Sheet="ARM"; Cell="J118"; Value=0}
    @{Sheet="ARM"; Cell="H122"; Value=4390.9414}
    @{Sheet="ARM"; Cell="I122"; Value=4400.759}
    @{Sheet="ARM"; Cell="J122"; Value=4334}
    @{Sheet="ARM"; Cell="K122"; Value=13202.277}
    @{Sheet="ARM"; Cell="L122"; Value=13002}
    @{Sheet="ARM"; Cell="M122"; Value=-10752.277}
    @{Sheet="ARM"; Cell="N122"; Value=-17902}
    @{Sheet="ARM"; Cell="L123"; Value=0}
    @{Sheet="ARM"; Cell="H123"; Value=0}
    @{Sheet="ARM"; Cell="J123"; Value=0}
    @{Sheet="ARM"; Cell="L131"; Value=0}
    @{Sheet="ARM"; Cell="H131"; Value=0}
    @{Sheet="ARM"; Cell="J131"; Value=0}
    @{Sheet="ARM"; Cell="H132"; Value=14134.863}
    @{Sheet="ARM"; Cell="I132"; Value=14198.45}
    @{Sheet="ARM"; Cell="K132"; Value=42595.35000000001}
    @{Sheet="ARM"; Cell="M132"; Value=-40065.35000000001}
    @{Sheet="BSM"; Cell="H3"; Value=2113.9}
    @{Sheet="BSM"; Cell="I3"; Value=1172.5264}
    @{Sheet="BSM"; Cell="K3"; Value=1172.5264}
    @{Sheet="BSM"; Cell="M3"; Value=-1058.5264}
    @{Sheet="BSM"; Cell="M75"; Value=-15225.667}
    @{Sheet="BSM"; Cell="H75"; Value=16161.667}
    @{Sheet="BSM"; Cell="I75"; Value=16161.667}
    @{Sheet="BSM"; Cell="J75"; Value=0}
    @{Sheet="BSM"; Cell="K75"; Value=16161.667}
    @{Sheet="BSM"; Cell="L75"; Value=0}
    @{Sheet="BSM"; Cell="M78"; Value=-43805.001}
    @{Sheet="BSM"; Cell="H78"; Value=16161.667}
    @{Sheet="BSM"; Cell="I78"; Value=16161.667}
    @{Sheet="BSM"; Cell="J78"; Value=0}
    @{Sheet="BSM"; Cell="K78"; Value=48485.001}
    @{Sheet="BSM"; Cell="L78"; Value=0}
    @{Sheet="BSM"; Cell="H99"; Value=5218.1787}
    @{Sheet="BSM"; Cell="I99"; Value=5728.6665}
    @{Sheet="BSM"; Cell="J99"; Value=2155.25}
    @{Sheet="BSM"; Cell="K99"; Value=5728.6665}
    @{Sheet="BSM"; Cell="L99"; Value=2155.25}
    @{Sheet="BSM"; Cell="M99"; Value=-4230.6665}
    @{Sheet="BSM"; Cell="N99"; Value=-5151.25}
    @{Sheet="CRP"; Cell="H31"; Value=2454.795}
    @{Sheet="CRP"; Cell="I31"; Value=1001.55}
    @{Sheet="CRP"; Cell="K31"; Value=1001.55}
    @{Sheet="CRP"; Cell="M31"; Value=-706.55}
    @{Sheet="CRP"; Cell="H34"; Value=2454.795}
    @{Sheet="CRP"; Cell="I34"; Value=1001.55}
    @{Sheet="CRP"; Cell="K34"; Value=1001.55}
    @{Sheet="CRP"; Cell="M34"; Value=-799.55}
    @{Sheet="CRP"; Cell="H86"; Value=76927350}
    @{Sheet="CRP"; Cell="I86"; Value=100004190}
    @{Sheet="CRP"; Cell="K86"; Value=100004190}
    @{Sheet="CRP"; Cell="M86"; Value=-100003067}
    @{Sheet="CRP"; Cell="H89"; Value=76927350}
    @{Sheet="CRP"; Cell="I89"; Value=100004190}
    @{Sheet="CRP"; Cell="K89"; Value=500020950}
    @{Sheet="CRP"; Cell="M89"; Value=-500015334}
    @{Sheet="CRP"; Cell="H99"; Value=2357.182}
    @{Sheet="CRP"; Cell="I99"; Value=1936.6666}
    @{Sheet="CRP"; Cell="J99"; Value=4249.5}
    @{Sheet="CRP"; Cell="K99"; Value=1936.6666}
    @{Sheet="CRP"; Cell="L99"; Value=4249.5}
    @{Sheet="CRP"; Cell="M99"; Value=-438.6666}
    @{Sheet="CRP"; Cell="N99"; Value=-7245.5}
    @{Sheet="CRP"; Cell="H105"; Value=1756.75}
    @{Sheet="CRP"; Cell="I105"; Value=1916.5385}
    @{Sheet="CRP"; Cell="K105"; Value=1916.5385}
    @{Sheet="CRP"; Cell="M105"; Value=-169.5385000000001}
    @{Sheet="CRP"; Cell="H107"; Value=756.3103599999999}
    @{Sheet="CRP"; Cell="I107"; Value=455.57895}
    @{Sheet="CRP"; Cell="K107"; Value=455.57895}
    @{Sheet="CRP"; Cell="M107"; Value=1464.42105}
    @{Sheet="CRP"; Cell="H122"; Value=2080.2727}
    @{Sheet="CRP"; Cell="I122"; Value=1965}
    @{Sheet="CRP"; Cell="J122"; Value=2218.6}
    @{Sheet="CRP"; Cell="K122"; Value=5895}
    @{Sheet="CRP"; Cell="L122"; Value=6655.799999999999}
    @{Sheet="CRP"; Cell="M122"; Value=-3445}
    @{Sheet="CRP"; Cell="N122"; Value=-11555.8}
    @{Sheet="CRP"; Cell="H126"; Value=2357.182}
    @{Sheet="CRP"; Cell="I126"; Value=1936.6666}
    @{Sheet="CRP"; Cell="J126"; Value=4249.5}
    @{Sheet="CRP"; Cell="K126"; Value=5809.9998}
    @{Sheet="CRP"; Cell="L126"; Value=12748.5}
    @{Sheet="CRP"; Cell="M126"; Value=-3339.9998}
    @{Sheet="CRP"; Cell="N126"; Value=-17688.5}
    @{Sheet="CRP"; Cell="H132"; Value=2926.1304}
    @{Sheet="CRP"; Cell="I132"; Value=2633.3809}
    @{Sheet="CRP"; Cell="K132"; Value=7900.1427}
    @{Sheet="CRP"; Cell="M132"; Value=-5370.1427}
    @{Sheet="CRP"; Cell="H134"; Value=2192.44}
    @{Sheet="CRP"; Cell="I134"; Value=1778.2333}
    @{Sheet="CRP"; Cell="K134"; Value=5334.699900000001}
    @{Sheet="CRP"; Cell="M134"; Value=-2799.699900000001}
    @{Sheet="CUL"; Cell="H2"; Value=540598.8}
    @{Sheet="CUL"; Cell="I2"; Value=740764.75}
    @{Sheet="CUL"; Cell="K2"; Value=4444588.5}
    @{Sheet="CUL"; Cell="M2"; Value=-4444475.5}
    @{Sheet="CUL"; Cell="H38"; Value=252}
    @{Sheet="CUL"; Cell="I38"; Value=152.28572}
    @{Sheet="CUL"; Cell="K38"; Value=456.85716}
    @{Sheet="CUL"; Cell="M38"; Value=-109.85716}
    @{Sheet="CUL"; Cell="H46"; Value=186852000}
    @{Sheet="CUL"; Cell="J46"; Value=336666850}
    @{Sheet="CUL"; Cell="L46"; Value=1010000550}
    @{Sheet="CUL"; Cell="N46"; Value=-1010000732}
    @{Sheet="CUL"; Cell="H63"; Value=12005.1}
    @{Sheet="CUL"; Cell="I63"; Value=11968.8}
    @{Sheet="CUL"; Cell="J63"; Value=12041.4}
    @{Sheet="CUL"; Cell="K63"; Value=35906.39999999999}
    @{Sheet="CUL"; Cell="L63"; Value=36124.2}
    @{Sheet="CUL"; Cell="M63"; Value=-35157.39999999999}
    @{Sheet="CUL"; Cell="N63"; Value=-37622.2}
    @{Sheet="CUL"; Cell="H66"; Value=12005.1}
    @{Sheet="CUL"; Cell="I66"; Value=11968.8}
    @{Sheet="CUL"; Cell="J66"; Value=12041.4}
    @{Sheet="CUL"; Cell="K66"; Value=107719.2}
    @{Sheet="CUL"; Cell="L66"; Value=108372.6}
    @{Sheet="CUL"; Cell="M66"; Value=-103975.2}
    @{Sheet="CUL"; Cell="N66"; Value=-115860.6}
    @{Sheet="GSM"; Cell="H102"; Value=4034.95}
    @{Sheet="GSM"; Cell="I102"; Value=2568.6316}
    @{Sheet="GSM"; Cell="K102"; Value=2568.6316}
    @{Sheet="GSM"; Cell="M102"; Value=-946.6316000000002}
    @{Sheet="GSM"; Cell="H122"; Value=5718.857}
    @{Sheet="GSM"; Cell="I122"; Value=7112.5}
    @{Sheet="GSM"; Cell="J122"; Value=2234.75}
    @{Sheet="GSM"; Cell="K122"; Value=21337.5}
    @{Sheet="GSM"; Cell="L122"; Value=6704.25}
    @{Sheet="GSM"; Cell="M122"; Value=-18887.5}
    @{Sheet="GSM"; Cell="N122"; Value=-11604.25}
    @{Sheet="GSM"; Cell="H126"; Value=13007341}
    @{Sheet="GSM"; Cell="I126"; Value=7812}
    @{Sheet="GSM"; Cell="K126"; Value=23436}
    @{Sheet="GSM"; Cell="M126"; Value=-20966}
    @{Sheet="GSM"; Cell="H132"; Value=3876}
    @{Sheet="GSM"; Cell="I132"; Value=3905.2646}
    @{Sheet="GSM"; Cell="J132"; Value=3828.6191}
    @{Sheet="GSM"; Cell="K132"; Value=11715.7938}
    @{Sheet="GSM"; Cell="L132"; Value=11485.8573}
    @{Sheet="GSM"; Cell="M132"; Value=-9185.793799999999}
    @{Sheet="GSM"; Cell="N132"; Value=-16545.8573}
    @{Sheet="LTW"; Cell="H40"; Value=7135}
    @{Sheet="LTW"; Cell="I40"; Value=7515.769}
    @{Sheet="LTW"; Cell="J40"; Value=5485}
    @{Sheet="LTW"; Cell="K40"; Value=7515.769}
    @{Sheet="LTW"; Cell="L40"; Value=5485}
    @{Sheet="LTW"; Cell="M40"; Value=-7379.769}
    @{Sheet="LTW"; Cell="N40"; Value=-5757}
    @{Sheet="LTW"; Cell="H122"; Value=8039.4614}
    @{Sheet="LTW"; Cell="I122"; Value=8074.75}
    @{Sheet="LTW"; Cell="J122"; Value=7983}
    @{Sheet="LTW"; Cell="K122"; Value=24224.25}
    @{Sheet="LTW"; Cell="L122"; Value=23949}
    @{Sheet="LTW"; Cell="M122"; Value=-21774.25}
    @{Sheet="LTW"; Cell="N122"; Value=-28849}
    @{Sheet="LTW"; Cell="H136"; Value=3149.6296}
    @{Sheet="LTW"; Cell="J136"; Value=3767.3572}
    @{Sheet="LTW"; Cell="L136"; Value=11302.0716}
    @{Sheet="LTW"; Cell="N136"; Value=-16402.0716}
    @{Sheet="WVR"; Cell="H122"; Value=4005.963}
    @{Sheet="WVR"; Cell="I122"; Value=3941.8696}
    @{Sheet="WVR"; Cell="K122"; Value=11825.6088}
    @{Sheet="WVR"; Cell="M122"; Value=-9375.6088}
    @{Sheet="WVR"; Cell="H132"; Value=5092.4087}
    @{Sheet="WVR"; Cell="I132"; Value=5151.291}
    @{Sheet="WVR"; Cell="J132"; Value=4890}
    @{Sheet="WVR"; Cell="K132"; Value=15453.873}
    @{Sheet="WVR"; Cell="L132"; Value=14670}
    @{Sheet="WVR"; Cell="M132"; Value=-12923.873}
    @{Sheet="WVR"; Cell="N132"; Value=-19730}
)

$deletions = @(
    @{Sheet="ARM"; Cell="N118"}
    @{Sheet="ARM"; Cell="N123"}
    @{Sheet="ARM"; Cell="N131"}
    @{Sheet="BSM"; Cell="N75"}
    @{Sheet="BSM"; Cell="N78"}
)

foreach ($c in $changes) {
    $ws = $wb.Worksheets.Item($c.Sheet)
    $ws.Range($c.Cell).Value = $c.Value
}

foreach ($d in $deletions) {
    $ws = $wb.Worksheets.Item($d.Sheet)
    $ws.Range($d.Cell).ClearContents()
}
